# The "2024" sheet tracks monthly SMS-derived transaction notes, with the
# most recent entry for each month always inserted at the top of that
# month's block (pushing older entries down). A new September entry was
# captured, so insert a fresh row above the first existing September row
# (row 35) and fill in its Details/Date pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$ws.Rows("35:35").Insert()

$ws.Range("R35").Value = "ift anbu tpar"
$ws.Range("S35").Value = "2024-09-09 11:27:52"
